# Hawaii Samples 11/19/2019 and 11/20/2019 (last HI Samples)
# Appends two new sample rows (69 and 70) to Sheet1, mirroring the existing
# table layout (Date | CRM value | Batch value | % off | Batch # | note).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 69: 11/19/2019 sample -------------------------------------------
# Copy row 68's date cell first so the new date cell inherits the existing
# short-date number format (style index 1) instead of Excel minting a brand
# new custom numFmt; then overwrite with the real value.
$ws.Cells.Item(68, 1).Copy($ws.Cells.Item(69, 1))
$ws.Cells.Item(69, 1).Value = 43788          # 11/19/2019
$ws.Cells.Item(69, 2).Value = 2208.5000007860199
$ws.Cells.Item(69, 3).Value = 2207.0300000000002
$ws.Range("D69").Formula = "=100*(B69-C69)/C69"
$ws.Cells.Item(69, 5).Value = 169
$ws.Cells.Item(69, 6).Value = "New crm opened 11/19/2019"

# --- Row 70: 11/20/2019 sample (last Hawaii sample) -----------------------
$ws.Cells.Item(68, 1).Copy($ws.Cells.Item(70, 1))
$ws.Cells.Item(70, 1).Value = 43789          # 11/20/2019
$ws.Cells.Item(70, 2).Value = 2203.6505910000001
$ws.Cells.Item(70, 3).Value = 2207.0300000000002
$ws.Range("D70").Formula = "=100*(B70-C70)/C70"
$ws.Cells.Item(70, 5).Value = 169
$ws.Cells.Item(70, 6).Value = "crm opened 11/19/2019"

# Match the author's final selection (cell below the newly-added data).
[void]$ws.Range("A71").Select()
